$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) and 1h volume % (column E) refresh.
# Some Price values are plain decimal-looking strings (e.g. "95.08"); Excel's
# normal smart-typing would silently coerce those into numeric cells, but the
# source data models them as text (mixed formats like "30.344.33" cannot be
# numeric). Force text entry via NumberFormat "@" for those cells, then restore
# the default "Normal" style so we do not leave a stray number format behind.
$textForceCells = @(
    "D5",
    "D6",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D24",
    "D25",
    "D26",
    "D27",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) updates
$ws.Range("D2").Value = "30.344.33"
$ws.Range("D3").Value = "1.941.13"
$ws.Range("D5").Value = "250.94"
$ws.Range("D6").Value = "0.7212"
$ws.Range("D10").Value = "0.07378"
$ws.Range("D11").Value = "0.8153"
$ws.Range("D12").Value = "0.08124"
$ws.Range("D13").Value = "1.937.47"
$ws.Range("D14").Value = "5.488"
$ws.Range("D15").Value = "95.08"
$ws.Range("D16").Value = "14.96"
$ws.Range("D17").Value = "30.356.88"
$ws.Range("D18").Value = "0.000008351"
$ws.Range("D19").Value = "252.48"
$ws.Range("D20").Value = "5.887"
$ws.Range("D21").Value = "2.194.79"
$ws.Range("D24").Value = "6.972"
$ws.Range("D25").Value = "9.850"
$ws.Range("D26").Value = "163.04"
$ws.Range("D27").Value = "2.407"
$ws.Range("D30").Value = "1.573"
$ws.Range("D32").Value = "4.460"
$ws.Range("D33").Value = "4.250"
$ws.Range("D34").Value = "0.05271"
$ws.Range("D35").Value = "1.305"
$ws.Range("D36").Value = "0.7560"
$ws.Range("D38").Value = "0.01994"
$ws.Range("D39").Value = "2.859"
$ws.Range("D40").Value = "81.47"
$ws.Range("D41").Value = "6.611"
$ws.Range("D42").Value = "0.4558"
$ws.Range("D43").Value = "2.042"
$ws.Range("D44").Value = "0.8477"
$ws.Range("D46").Value = "102.58"
$ws.Range("D47").Value = "9.845"
$ws.Range("D48").Value = "7.509"
$ws.Range("D49").Value = "36.89"
$ws.Range("D50").Value = "0.4193"
$ws.Range("D51").Value = "1.510"

# Volume(1h) % (column E) updates
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  -9.08%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  +5.54%  "
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("E18").Value = "  +5.21%  "
$ws.Range("E19").Value = "  -7.49%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  -12.41%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  +7.37%  "
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  -1.93%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
